# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for rows that changed per the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.237.81"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.485.53"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.23"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.85"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.484.09"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "8.08"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("E11").Value = "  -5.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.411"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.073.87"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.33"
$ws.Range("E14").Value = "  -3.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000202"
$ws.Range("E15").Value = "  -5.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.479.13"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.204.96"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.83"
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.14"
$ws.Range("E20").Value = "  -3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.79"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.98"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.82"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000116"
$ws.Range("E26").Value = "  -3.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  -5.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.95"
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.166"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("E32").Value = "  -8.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.06"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.467.77"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.73"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("E37").Value = "  -6.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.71"
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "170.30"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0860"
$ws.Range("E41").Value = "  -3.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.15"
$ws.Range("E42").Value = "  -5.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.880"
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.91"
$ws.Range("E44").Value = "  -8.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "45.39"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.12"
$ws.Range("E46").Value = "  -9.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.13"
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.937"
$ws.Range("E50").Value = "  -3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.237"
$ws.Range("E51").Value = "  -4.23%  "
